$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 14:17"
$ws.Range("B4").Value = 3356242
$ws.Range("C4").Value = 596
$ws.Range("E4").Value = 1728126
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 137414
$ws.Range("B6").Value = 854480
$ws.Range("C6").Value = 4122
$ws.Range("D6").Value = 537599
$ws.Range("E6").Value = 294163
$ws.Range("G6").Value = 31
$ws.Range("H6").Value = 22718
$ws.Range("B14").Value = 257303
$ws.Range("C14").Value = 2186
$ws.Range("D14").Value = 219993
$ws.Range("E14").Value = 24481
$ws.Range("G14").Value = 194
$ws.Range("H14").Value = 12829
$ws.Range("B32").Value = 64932
$ws.Range("C32").Value = 165
$ws.Range("D32").Value = 55380
$ws.Range("E32").Value = 9088
$ws.Range("G32").Value = 5
$ws.Range("H32").Value = 464
$ws.Range("A36").Value = "Kuwait"
$ws.Range("B36").Value = 54894
$ws.Range("C36").Value = 836
$ws.Range("D36").Value = 44610
$ws.Range("E36").Value = 9894
$ws.Range("G36").Value = 4
$ws.Range("H36").Value = 390
$ws.Range("A37").Value = "Emiratos Arabes Unidos"
$ws.Range("B37").Value = 54453
$ws.Range("D37").Value = 44648
$ws.Range("E37").Value = 9474
$ws.Range("H37").Value = 331
$ws.Range("A38").Value = "Filipinas"
$ws.Range("B38").Value = 54222
$ws.Range("D38").Value = 14037
$ws.Range("E38").Value = 38813
$ws.Range("H38").Value = 1372
$ws.Range("B49").Value = 32883
$ws.Range("C49").Value = 66
$ws.Range("E49").Value = 1415
$ws.Range("A50").Value = "Rumania"
$ws.Range("B50").Value = 32535
$ws.Range("C50").Value = 456
$ws.Range("D50").Value = 21545
$ws.Range("E50").Value = 9106
$ws.Range("G50").Value = 13
$ws.Range("H50").Value = 1884
$ws.Range("A51").Value = "Barein"
$ws.Range("B51").Value = 32470
$ws.Range("D51").Value = 27828
$ws.Range("E51").Value = 4538
$ws.Range("H51").Value = 104
$ws.Range("B64").Value = 16801
$ws.Range("C64").Value = 82
$ws.Range("D64").Value = 8589
$ws.Range("E64").Value = 8174
$ws.Range("B70").Value = 12855
$ws.Range("C70").Value = 342
$ws.Range("E70").Value = 5073
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 59
$ws.Range("B79").Value = 8718
$ws.Range("C79").Value = 14
$ws.Range("D79").Value = 8519
$ws.Range("E79").Value = 77
$ws.Range("B80").Value = 8135
$ws.Range("C80").Value = 121
$ws.Range("D80").Value = 5446
$ws.Range("E80").Value = 2541
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 148
$ws.Range("A81").Value = "Consejo Danes para los Refugiados"
$ws.Range("B81").Value = 8033
$ws.Range("C81").Value = 62
$ws.Range("D81").Value = 3615
$ws.Range("E81").Value = 4229
$ws.Range("H81").Value = 189
$ws.Range("A82").Value = "Republica de Macedonia"
$ws.Range("B82").Value = 7975
$ws.Range("D82").Value = 4080
$ws.Range("E82").Value = 3519
$ws.Range("H82").Value = 376
$ws.Range("B87").Value = 6877
$ws.Range("C87").Value = 158
$ws.Range("D87").Value = 3115
$ws.Range("E87").Value = 3541
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 221
$ws.Range("A96").Value = "Madagascar"
$ws.Range("B96").Value = 4867
$ws.Range("C96").Value = 289
$ws.Range("D96").Value = 2378
$ws.Range("E96").Value = 2454
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 35
$ws.Range("A97").Value = "Luxemburgo"
$ws.Range("B97").Value = 4842
$ws.Range("D97").Value = 4086
$ws.Range("E97").Value = 646
$ws.Range("H97").Value = 110
$ws.Range("B101").Value = 3722
$ws.Range("C101").Value = 50
$ws.Range("D101").Value = 2486
$ws.Range("E101").Value = 1117
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 119
$ws.Range("B102").Value = 3454
$ws.Range("C102").Value = 83
$ws.Range("D102").Value = 1946
$ws.Range("E102").Value = 1415
$ws.Range("G102").Value = 4
$ws.Range("H102").Value = 93
$ws.Range("B110").Value = 2605
$ws.Range("C110").Value = 94
$ws.Range("D110").Value = 1981
$ws.Range("E110").Value = 613
$ws.Range("B121").Value = 1869
$ws.Range("C121").Value = 4
$ws.Range("D121").Value = 1571
$ws.Range("E121").Value = 219
$ws.Range("B133").Value = 1263
$ws.Range("C133").Value = 18
$ws.Range("E133").Value = 137
$ws.Range("A148").Value = "Namibia"
$ws.Range("B148").Value = 785
$ws.Range("C148").Value = 72
$ws.Range("D148").Value = 26
$ws.Range("E148").Value = 758
$ws.Range("H148").Value = 1
$ws.Range("A149").Value = "Jamaica"
$ws.Range("B149").Value = 758
$ws.Range("C149").Value = 5
$ws.Range("D149").Value = 615
$ws.Range("E149").Value = 133
$ws.Range("H149").Value = 10
$ws.Range("A150").Value = "Surinam"
$ws.Range("B150").Value = 741
$ws.Range("D150").Value = 495
$ws.Range("E150").Value = 228
$ws.Range("H150").Value = 18
$ws.Range("A151").Value = "Santo Tome y Principe"
$ws.Range("B151").Value = 727
$ws.Range("D151").Value = 284
$ws.Range("E151").Value = 429
$ws.Range("H151").Value = 14
$ws.Range("A152").Value = "Crucero"
$ws.Range("B152").Value = 712
$ws.Range("D152").Value = 651
$ws.Range("E152").Value = 48
$ws.Range("H152").Value = 13
$ws.Range("A153").Value = "Togo"
$ws.Range("B153").Value = 710
$ws.Range("D153").Value = 494
$ws.Range("E153").Value = 201
$ws.Range("H153").Value = 15
$ws.Range("A154").Value = "San Marino"
$ws.Range("B154").Value = 699
$ws.Range("D154").Value = 656
$ws.Range("E154").Value = 1
$ws.Range("H154").Value = 42
$ws.Range("A155").Value = "Malta"
$ws.Range("B155").Value = 674
$ws.Range("D155").Value = 660
$ws.Range("E155").Value = 5
$ws.Range("H155").Value = 9
$ws.Range("B161").Value = 372
$ws.Range("C161").Value = 2
$ws.Range("E161").Value = 22
